# Change the computations of the KPIs
# Updates StartingInventories/SetupCosts on Productdata, Capacity values,
# and the diagonal ProcessingTime matrix.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Productdata sheet: columns C (StartingInventories) and E (SetupCosts)
# ---------------------------------------------------------------------------
$wsProd = $wb.Worksheets.Item("Productdata")

$prodC = @{
    2  = 1051
    3  = 1891
    4  = 856
    5  = 691
    6  = 1
    7  = 1
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 1
    15 = 0
    16 = 0
    17 = 0
    18 = 0
}

$prodE = @{
    2  = 86.625
    3  = 232.6401
    4  = 106.7553
    5  = 87.92117999999999
    6  = 136.08
    7  = 182.9475
    8  = 224.1351
    9  = 181.3212
    10 = 62.37
    11 = 80.73000000000002
    12 = 53.82000000000001
    13 = 85.63500000000001
    14 = 154.575
    15 = 136.08
    16 = 85.22010000000002
    17 = 69.6177
    18 = 182.9475
}

foreach ($row in $prodC.Keys) {
    $wsProd.Cells.Item($row, 3).Value = $prodC[$row]
}

foreach ($row in $prodE.Keys) {
    $wsProd.Cells.Item($row, 5).Value = $prodE[$row]
}

# ---------------------------------------------------------------------------
# Capacity sheet: column B
# ---------------------------------------------------------------------------
$wsCap = $wb.Worksheets.Item("Capacity")

$capB = @{
    2  = 700
    3  = 6300
    4  = 2850
    5  = 460
    6  = 2520
    7  = 8650
    8  = 3780
    9  = 2060
    10 = 2520
    11 = 2990
    12 = 2990
    13 = 5190
    14 = 9160
    15 = 6300
    16 = 1260
    17 = 4120
    18 = 8650
}

foreach ($row in $capB.Keys) {
    $wsCap.Cells.Item($row, 2).Value = $capB[$row]
}

# ---------------------------------------------------------------------------
# ProcessingTime sheet: diagonal values (row r, col r) for rows 2..18
# ---------------------------------------------------------------------------
$wsProc = $wb.Worksheets.Item("ProcessingTime")

$procDiag = @{
    2  = 1
    3  = 5
    5  = 1
    6  = 2
    7  = 5
    8  = 3
    9  = 2
    10 = 2
    11 = 1
    12 = 1
    13 = 3
    14 = 4
    15 = 5
    16 = 1
    17 = 4
    18 = 5
}

foreach ($row in $procDiag.Keys) {
    $wsProc.Cells.Item($row, $row).Value = $procDiag[$row]
}
